# toerstraaling.xlsx - "Fixed før og efterbehandling, samt screenshot funktion"
#
# Cell D4 held a formula that pulled from an external, linked workbook
# ( [1]PL_EURO / [1]PL_DKK_20161101 ). Editing that cell in the Excel UI
# trips Excel's "this cell has an external reference that can't be shown
# or edited" guard: the cell's formula is replaced by that message (as a
# literal string) and the external link is severed for it, which in turn
# breaks the shared formulas in B4/C4 (they were sharing the formula
# definition anchored at this column) -> they turn into #VALUE! errors.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = "This cell has an external reference that can't be shown or edited. Editing this cell will remove the external reference."

# Selection left on F2 when the file was saved.
$ws.Range("F2").Select()
